$d = $word.ActiveDocument

# 1. Title: merge "Adversary" + " Attribution Report" -> "Threat Attribution Report"
$d.Content.Find.Execute("Adversary Attribution Report", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Threat Attribution Report", 2)

# 2. Remove the _GoBack bookmark (bookmarkStart/bookmarkEnd) after "Detected TTPs..." paragraph
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3. Heading "2. Overlap " + "Table" -> "2. Overlap Table"
$d.Content.Find.Execute("2. Overlap Table", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2. Overlap Table", 2)

# 4. "[Reflections and insights for future read" + "iness.]" -> single run
$d.Content.Find.Execute("[Reflections and insights for future readiness.]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[Reflections and insights for future readiness.]", 2)
